# edit.ps1
# Applies the changes described by the commit "update report, log result":
#  - Detail sheet: update G19 (actual result text) for TC18
#  - Cart sheet: update G11 (TC10), G12/H12 (TC11), G17/H17 (TC16)
#  - Payment sheet: update G17/H17 (TC16), G18/H18 (TC17)
#  - Add two new sheets "Demo" and "Demo2" at the end with test-case data,
#    and make "Demo" the active sheet/tab.

$wb = $excel.ActiveWorkbook

# --- Detail sheet: TC18 (row 19) actual result ---
$wsDetail = $wb.Worksheets.Item("Detail")
$wsDetail.Cells.Item(19, 7).Value = "Thêm sản phẩm vào giỏ hàng thành công."

# --- Cart sheet updates ---
$wsCart = $wb.Worksheets.Item("Cart")
$wsCart.Cells.Item(11, 7).Value = "Tổng tiền sản phẩm là là hợp lệ"
$wsCart.Cells.Item(12, 7).Value = '"Tổng tiền giỏ hàng là hợp lệ, bằng tổng tiền các sản phẩm"'
$wsCart.Cells.Item(12, 8).Value = "Pass"
$wsCart.Cells.Item(17, 7).Value = '"Chuyển hướng đến trang thanh toán"'
$wsCart.Cells.Item(17, 8).Value = "Pass"

# --- Payment sheet updates ---
$wsPayment = $wb.Worksheets.Item("Payment")
$wsPayment.Cells.Item(17, 7).Value = "Chuyển hướng đến trang thanh toán MOMO"
$wsPayment.Cells.Item(17, 8).Value = "Pass"
$wsPayment.Cells.Item(18, 7).Value = "Không tìm thấy thông báo nào"
$wsPayment.Cells.Item(18, 8).Value = "Fail"

# --- Add "Demo" sheet at the end (copy of Search-style test cases) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws1 = $wb.Worksheets.Add($null, $lastSheet)
$ws1.Name = "Demo"

$ws1.Cells.Item(1, 1).Value = 'ID'
$ws1.Cells.Item(1, 2).Value = 'Chức năng'
$ws1.Cells.Item(1, 3).Value = 'Loại test case'
$ws1.Cells.Item(1, 4).Value = 'Mô tả'
$ws1.Cells.Item(1, 5).Value = 'Dữ liệu test'
$ws1.Cells.Item(1, 6).Value = 'Kỳ vọng'
$ws1.Cells.Item(1, 7).Value = 'Thực tế'
$ws1.Cells.Item(1, 8).Value = 'Kết quả'
$ws1.Cells.Item(2, 1).Value = 'TC01'
$ws1.Cells.Item(2, 2).Value = 'Tìm kiếm sản phẩm'
$ws1.Cells.Item(2, 3).Value = 'UI'
$ws1.Cells.Item(2, 4).Value = 'Kiểm tra sự hiển thị của ô tìm kiếm'
$ws1.Cells.Item(2, 6).Value = '"Ô tìm kiếm hiển thị, placeholder là ''Tìm kiếm sản phẩm''"'
$ws1.Cells.Item(2, 7).Value = 'ô tìm kiếm hiển thị placeholder là tìm kiếm sản phẩm'
$ws1.Cells.Item(2, 8).Value = 'Pass'
$ws1.Cells.Item(3, 1).Value = 'TC02'
$ws1.Cells.Item(3, 2).Value = 'Tìm kiếm sản phẩm'
$ws1.Cells.Item(3, 3).Value = 'UI'
$ws1.Cells.Item(3, 4).Value = 'Kiểm tra sự hiển thị của button tìm kiếm'
$ws1.Cells.Item(3, 6).Value = '"Button tìm kiếm hiển thị"'
$ws1.Cells.Item(3, 7).Value = 'button tìm kiếm hiển thị'
$ws1.Cells.Item(3, 8).Value = 'Pass'
$ws1.Cells.Item(4, 1).Value = 'TC03'
$ws1.Cells.Item(4, 2).Value = 'Tìm kiếm sản phẩm'
$ws1.Cells.Item(4, 3).Value = 'UI'
$ws1.Cells.Item(4, 4).Value = 'Kiểm tra placeholder của ô tìm kiếm'
$ws1.Cells.Item(4, 6).Value = '"Placeholder trong ô tìm kiếm là ''Tìm kiếm sản phẩm''"'
$ws1.Cells.Item(4, 7).Value = 'placeholder trong ô tìm kiếm là tìm kiếm sản phẩm'
$ws1.Cells.Item(4, 8).Value = 'Pass'
$ws1.Cells.Item(5, 1).Value = 'TC04'
$ws1.Cells.Item(5, 2).Value = 'Tìm kiếm sản phẩm'
$ws1.Cells.Item(5, 3).Value = 'UI'
$ws1.Cells.Item(5, 4).Value = 'Kiểm tra ngôn ngữ hiển thị của placeholder'
$ws1.Cells.Item(5, 6).Value = '"Placeholder hiển thị bằng tiếng Việt"'
$ws1.Cells.Item(5, 7).Value = 'placeholder hiển thị bằng tiếng việt'
$ws1.Cells.Item(5, 8).Value = 'Pass'
$ws1.Cells.Item(6, 1).Value = 'TC05'
$ws1.Cells.Item(6, 2).Value = 'Tìm kiếm sản phẩm'
$ws1.Cells.Item(6, 3).Value = 'Chức năng'
$ws1.Cells.Item(6, 4).Value = 'Tìm kiếm với từ khóa có trong danh sách sản phẩm (Life)'
$ws1.Cells.Item(6, 5).Value = 'Từ khóa: Life'
$ws1.Cells.Item(6, 6).Value = '"Danh sách sản phẩm chứa ''Life'' hiển thị"'
$ws1.Cells.Item(6, 7).Value = '"Danh sách sản phẩm chứa ''Life'' hiển thị"'
$ws1.Cells.Item(6, 8).Value = 'Pass'
$ws1.Cells.Item(7, 1).Value = 'TC06'
$ws1.Cells.Item(7, 2).Value = 'Tìm kiếm sản phẩm'
$ws1.Cells.Item(7, 3).Value = 'Chức năng'
$ws1.Cells.Item(7, 4).Value = 'Tìm kiếm với từ khóa không có trong danh sách sản phẩm (zzxy)'
$ws1.Cells.Item(7, 5).Value = 'Từ khóa: zzxy'
$ws1.Cells.Item(7, 6).Value = '"Danh sách sản phẩm không hiển thị, button ''Xem thêm'' hiển thị"'
$ws1.Cells.Item(7, 7).Value = 'Không tìm được sản phẩm'
$ws1.Cells.Item(7, 8).Value = 'Pass'
$ws1.Cells.Item(8, 1).Value = 'TC07'
$ws1.Cells.Item(8, 2).Value = 'Tìm kiếm sản phẩm'
$ws1.Cells.Item(8, 3).Value = 'Chức năng'
$ws1.Cells.Item(8, 4).Value = 'Tìm kiếm với từ khóa có khoảng trắng ở đầu (  Life)'
$ws1.Cells.Item(8, 5).Value = 'Từ khóa:   Life'
$ws1.Cells.Item(8, 6).Value = '"Danh sách sản phẩm chứa ''Life'' hiển thị"'
$ws1.Cells.Item(8, 7).Value = '"Danh sách sản phẩm chứa ''Life'' hiển thị"'
$ws1.Cells.Item(8, 8).Value = 'Pass'
$ws1.Cells.Item(9, 1).Value = 'TC08'
$ws1.Cells.Item(9, 2).Value = 'Tìm kiếm sản phẩm'
$ws1.Cells.Item(9, 3).Value = 'Chức năng'
$ws1.Cells.Item(9, 4).Value = 'Tìm kiếm với từ khóa có khoảng trắng ở cuối (Life  )'
$ws1.Cells.Item(9, 5).Value = 'Từ khóa: Life'
$ws1.Cells.Item(9, 6).Value = '"Danh sách sản phẩm chứa ''Life'' hiển thị"'
$ws1.Cells.Item(9, 7).Value = '"Danh sách sản phẩm chứa ''Life'' hiển thị"'
$ws1.Cells.Item(9, 8).Value = 'Pass'
$ws1.Cells.Item(10, 1).Value = 'TC09'
$ws1.Cells.Item(10, 2).Value = 'Tìm kiếm sản phẩm'
$ws1.Cells.Item(10, 3).Value = 'Chức năng'
$ws1.Cells.Item(10, 4).Value = 'Tìm kiếm với từ khóa có cả khoảng trắng đầu và cuối (  Life  )'
$ws1.Cells.Item(10, 5).Value = 'Từ khóa:   Life'
$ws1.Cells.Item(10, 6).Value = '"Danh sách sản phẩm chứa ''Life'' hiển thị"'
$ws1.Cells.Item(10, 7).Value = '"Danh sách sản phẩm chứa ''Life'' hiển thị"'
$ws1.Cells.Item(10, 8).Value = 'Pass'
$ws1.Cells.Item(11, 1).Value = 'TC10'
$ws1.Cells.Item(11, 2).Value = 'Tìm kiếm sản phẩm'
$ws1.Cells.Item(11, 3).Value = 'Chức năng'
$ws1.Cells.Item(11, 4).Value = 'Tìm kiếm với từ khóa là chữ thường (life)'
$ws1.Cells.Item(11, 5).Value = 'Từ khóa: life'
$ws1.Cells.Item(11, 6).Value = '"Danh sách sản phẩm chứa ''Life'' hiển thị"'
$ws1.Cells.Item(11, 7).Value = '"Danh sách sản phẩm chứa ''Life'' hiển thị"'
$ws1.Cells.Item(11, 8).Value = 'Pass'
$ws1.Cells.Item(12, 1).Value = 'TC11'
$ws1.Cells.Item(12, 2).Value = 'Tìm kiếm sản phẩm'
$ws1.Cells.Item(12, 3).Value = 'Chức năng'
$ws1.Cells.Item(12, 4).Value = 'Tìm kiếm với từ khóa là chữ hoa (LIFE)'
$ws1.Cells.Item(12, 5).Value = 'Từ khóa: LIFE'
$ws1.Cells.Item(12, 6).Value = '"Danh sách sản phẩm chứa ''Life'' hiển thị"'
$ws1.Cells.Item(12, 7).Value = '"Danh sách sản phẩm chứa ''Life'' hiển thị"'
$ws1.Cells.Item(12, 8).Value = 'Pass'
$ws1.Cells.Item(13, 1).Value = 'TC12'
$ws1.Cells.Item(13, 2).Value = 'Tìm kiếm sản phẩm'
$ws1.Cells.Item(13, 3).Value = 'Chức năng'
$ws1.Cells.Item(13, 4).Value = 'Tìm kiếm bằng cách nhấn Enter'
$ws1.Cells.Item(13, 5).Value = 'Từ khóa: Life'
$ws1.Cells.Item(13, 6).Value = '"Danh sách sản phẩm chứa ''Life'' hiển thị"'
$ws1.Cells.Item(13, 7).Value = '"Danh sách sản phẩm chứa ''Life'' hiển thị"'
$ws1.Cells.Item(13, 8).Value = 'Pass'
$ws1.Cells.Item(14, 1).Value = 'TC13'
$ws1.Cells.Item(14, 2).Value = 'Tìm kiếm sản phẩm'
$ws1.Cells.Item(14, 3).Value = 'Chức năng'
$ws1.Cells.Item(14, 4).Value = 'Tìm kiếm bằng cách nhấn nút tìm kiếm'
$ws1.Cells.Item(14, 5).Value = 'Từ khóa: Life'
$ws1.Cells.Item(14, 6).Value = '"Danh sách sản phẩm chứa ''Life'' hiển thị"'
$ws1.Cells.Item(14, 7).Value = '"Danh sách sản phẩm chứa ''Life'' hiển thị"'
$ws1.Cells.Item(14, 8).Value = 'Pass'
$ws1.Cells.Item(15, 1).Value = 'TC14'
$ws1.Cells.Item(15, 2).Value = 'Tìm kiếm sản phẩm'
$ws1.Cells.Item(15, 3).Value = 'Chức năng'
$ws1.Cells.Item(15, 4).Value = 'Tìm kiếm với từ khóa là ký tự đặc biệt (!@#)'
$ws1.Cells.Item(15, 5).Value = 'Từ khóa: !@#'
$ws1.Cells.Item(15, 6).Value = '"Danh sách sản phẩm không hiển thị, button ''Xem thêm'' hiển thị"'
$ws1.Cells.Item(15, 7).Value = 'Hệ thống cho phép tìm kiếm với ký tự đặc biệt, có kết quả trả về'
$ws1.Cells.Item(15, 8).Value = 'Pass'
$ws1.Cells.Item(16, 1).Value = 'TC15'
$ws1.Cells.Item(16, 2).Value = 'Tìm kiếm sản phẩm'
$ws1.Cells.Item(16, 3).Value = 'Chức năng'
$ws1.Cells.Item(16, 4).Value = 'Tìm kiếm với từ khóa là số (12345)'
$ws1.Cells.Item(16, 5).Value = 'Từ khóa: 12345'
$ws1.Cells.Item(16, 6).Value = '"Danh sách sản phẩm không hiển thị, button ''Xem thêm'' hiển thị"'
$ws1.Cells.Item(16, 7).Value = 'Hệ thống cho phép tìm kiếm với ký tự số, có kết quả trả về'
$ws1.Cells.Item(16, 8).Value = 'Pass'
$ws1.Cells.Item(17, 1).Value = 'TC16'
$ws1.Cells.Item(17, 2).Value = 'Tìm kiếm sản phẩm'
$ws1.Cells.Item(17, 3).Value = 'Chức năng'
$ws1.Cells.Item(17, 4).Value = 'Tìm kiếm với từ khóa là ký tự in hoa và số (LIFE123)'
$ws1.Cells.Item(17, 5).Value = 'Từ khóa: LIFE123'
$ws1.Cells.Item(17, 6).Value = '"Danh sách sản phẩm không hiển thị, button ''Xem thêm'' hiển thị"'
$ws1.Cells.Item(17, 7).Value = 'Không tìm được sản phẩm'
$ws1.Cells.Item(17, 8).Value = 'Pass'
$ws1.Cells.Item(18, 1).Value = 'TC17'
$ws1.Cells.Item(18, 2).Value = 'Tìm kiếm sản phẩm'
$ws1.Cells.Item(18, 3).Value = 'Chức năng'
$ws1.Cells.Item(18, 4).Value = 'Tìm kiếm với từ khóa trống'
$ws1.Cells.Item(18, 5).Value = 'Từ khóa:'
$ws1.Cells.Item(18, 6).Value = '"Danh sách sản phẩm không hiển thị, button ''Xem thêm'' hiển thị"'
$ws1.Cells.Item(18, 7).Value = '"Danh sách sản phẩm không hiển thị, button ''Xem thêm'' hiển thị"'
$ws1.Cells.Item(18, 8).Value = 'Pass'

# --- Add "Demo2" sheet at the end (copy of Cart-style test cases) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "Demo2"

$ws2.Cells.Item(1, 1).Value = 'ID'
$ws2.Cells.Item(1, 2).Value = 'Chức năng'
$ws2.Cells.Item(1, 3).Value = 'Loại test case'
$ws2.Cells.Item(1, 4).Value = 'Mô tả'
$ws2.Cells.Item(1, 5).Value = 'Dữ liệu test'
$ws2.Cells.Item(1, 6).Value = 'Kỳ vọng'
$ws2.Cells.Item(1, 7).Value = 'Thực tế'
$ws2.Cells.Item(1, 8).Value = 'Kết quả'
$ws2.Cells.Item(2, 1).Value = 'TC01'
$ws2.Cells.Item(2, 2).Value = 'Giỏ hàng'
$ws2.Cells.Item(2, 3).Value = 'UI'
$ws2.Cells.Item(2, 4).Value = 'Kiểm tra hiển thị các cột'
$ws2.Cells.Item(2, 6).Value = 'Hiển thị đầy đủ các cột: Sản Phẩm, Đơn giá, Số lượng, Tổng'
$ws2.Cells.Item(2, 7).Value = 'Không hỗ trợ mô tả'
$ws2.Cells.Item(2, 8).Value = 'Skip'
$ws2.Cells.Item(3, 1).Value = 'TC02'
$ws2.Cells.Item(3, 2).Value = 'Giỏ hàng'
$ws2.Cells.Item(3, 3).Value = 'UI'
$ws2.Cells.Item(3, 4).Value = 'Kiểm tra ngôn ngữ hiển thị'
$ws2.Cells.Item(3, 6).Value = 'Ngôn ngữ hiển thị là tiếng Việt'
$ws2.Cells.Item(3, 7).Value = 'ngôn ngữ hiển thị là tiếng việt'
$ws2.Cells.Item(3, 8).Value = 'Pass'
$ws2.Cells.Item(4, 1).Value = 'TC03'
$ws2.Cells.Item(4, 2).Value = 'Giỏ hàng'
$ws2.Cells.Item(4, 3).Value = 'UI'
$ws2.Cells.Item(4, 4).Value = 'Kiểm tra hiển thị tổng tiền'
$ws2.Cells.Item(4, 6).Value = 'Hiển thị dòng tổng tiền'
$ws2.Cells.Item(4, 7).Value = 'tổng tiền của giỏ hàng hiển thị đúng định dạng tiền tệ (vnđ) và giá trị chính xác (tổng của tất cả các tổng tiền sản phẩm)'
$ws2.Cells.Item(4, 8).Value = 'Fail'
$ws2.Cells.Item(5, 1).Value = 'TC04'
$ws2.Cells.Item(5, 2).Value = 'Giỏ hàng'
$ws2.Cells.Item(5, 3).Value = 'UI'
$ws2.Cells.Item(5, 4).Value = 'Kiểm tra hiển thị button "Thanh toán"'
$ws2.Cells.Item(5, 6).Value = 'Hiển thị button "Thanh toán"'
$ws2.Cells.Item(5, 7).Value = 'hiển thị button thanh toán'
$ws2.Cells.Item(5, 8).Value = 'Pass'
$ws2.Cells.Item(6, 1).Value = 'TC05'
$ws2.Cells.Item(6, 2).Value = 'Giỏ hàng'
$ws2.Cells.Item(6, 3).Value = 'UI'
$ws2.Cells.Item(6, 4).Value = 'Kiểm tra hiển thị liên kết "Tiếp tục mua sắm"'
$ws2.Cells.Item(6, 6).Value = 'Hiển thị liên kết "Tiếp tục mua sắm"'
$ws2.Cells.Item(6, 7).Value = 'hiển thị liên kết tiếp tục mua sắm'
$ws2.Cells.Item(6, 8).Value = 'Pass'
$ws2.Cells.Item(7, 1).Value = 'TC06'
$ws2.Cells.Item(7, 2).Value = 'Giỏ hàng'
$ws2.Cells.Item(7, 3).Value = 'UI'
$ws2.Cells.Item(7, 4).Value = 'Kiểm tra hiển thị button xóa'
$ws2.Cells.Item(7, 6).Value = 'Hiển thị button xóa cho từng sản phẩm'
$ws2.Cells.Item(7, 7).Value = 'hiển thị button xóa cho từng sản phẩm'
$ws2.Cells.Item(7, 8).Value = 'Pass'
$ws2.Cells.Item(8, 1).Value = 'TC07'
$ws2.Cells.Item(8, 2).Value = 'Giỏ hàng'
$ws2.Cells.Item(8, 3).Value = 'Chức năng'
$ws2.Cells.Item(8, 4).Value = 'Kiểm tra tính toán tổng tiền khi số lượng sản phẩm là 2'
$ws2.Cells.Item(8, 5).Value = 'Số lượng sản phẩm 1: 2, Số lượng sản phẩm 2: 2'
$ws2.Cells.Item(8, 6).Value = 'Tổng tiền hiển thị chính xác (tổng đơn giá * số lượng)'
$ws2.Cells.Item(8, 7).Value = 'Tổng tiền giỏ hàng là hợp lệ'
$ws2.Cells.Item(8, 8).Value = 'Fail'
$ws2.Cells.Item(9, 1).Value = 'TC08'
$ws2.Cells.Item(9, 2).Value = 'Giỏ hàng'
$ws2.Cells.Item(9, 3).Value = 'Chức năng'
$ws2.Cells.Item(9, 4).Value = 'Kiểm tra tính toán tổng của 1 sản phẩm'
$ws2.Cells.Item(9, 5).Value = 'Số lượng sản phẩm 1: 1'
$ws2.Cells.Item(9, 6).Value = 'Tổng hiển thị chính xác bằng đơn giá'
$ws2.Cells.Item(9, 7).Value = 'Số lượng sản phẩm hợp lệ'
$ws2.Cells.Item(9, 8).Value = 'Fail'
$ws2.Cells.Item(10, 1).Value = 'TC09'
$ws2.Cells.Item(10, 2).Value = 'Giỏ hàng'
$ws2.Cells.Item(10, 3).Value = 'Chức năng'
$ws2.Cells.Item(10, 4).Value = 'Kiểm tra thay đổi số lượng sản phẩm hợp lệ'
$ws2.Cells.Item(10, 5).Value = 'Số lượng sản phẩm 1: 2'
$ws2.Cells.Item(10, 6).Value = 'Số lượng sản phẩm được thay đổi thành công và tổng tiền được cập nhật tương ứng'
$ws2.Cells.Item(10, 7).Value = 'Tổng tiền giỏ hàng là hợp lệ'
$ws2.Cells.Item(10, 8).Value = 'Fail'
$ws2.Cells.Item(11, 1).Value = 'TC10'
$ws2.Cells.Item(11, 2).Value = 'Giỏ hàng'
$ws2.Cells.Item(11, 3).Value = 'Chức năng'
$ws2.Cells.Item(11, 4).Value = 'Kiểm tra thay đổi số lượng sản phẩm nhỏ hơn 1'
$ws2.Cells.Item(11, 5).Value = 'Số lượng sản phẩm 1: 0'
$ws2.Cells.Item(11, 6).Value = 'Hiển thị thông báo lỗi và số lượng không thay đổi'
$ws2.Cells.Item(11, 7).Value = 'vẫn giảm về 0'
$ws2.Cells.Item(11, 8).Value = 'Fail'
$ws2.Cells.Item(12, 1).Value = 'TC11'
$ws2.Cells.Item(12, 2).Value = 'Giỏ hàng'
$ws2.Cells.Item(12, 3).Value = 'Chức năng'
$ws2.Cells.Item(12, 4).Value = 'Kiểm tra thay đổi số lượng sản phẩm bằng 1'
$ws2.Cells.Item(12, 5).Value = 'Số lượng sản phẩm 1: 1'
$ws2.Cells.Item(12, 6).Value = 'Số lượng sản phẩm được thay đổi thành công'
$ws2.Cells.Item(12, 7).Value = 'Số lượng sản phẩm hợp lệ'
$ws2.Cells.Item(12, 8).Value = 'Fail'
$ws2.Cells.Item(13, 1).Value = 'TC12'
$ws2.Cells.Item(13, 2).Value = 'Giỏ hàng'
$ws2.Cells.Item(13, 3).Value = 'Chức năng'
$ws2.Cells.Item(13, 4).Value = 'Kiểm tra thay đổi số lượng sản phẩm lớn hơn số lượng tồn kho'
$ws2.Cells.Item(13, 5).Value = 'Số lượng sản phẩm 1: 453'
$ws2.Cells.Item(13, 6).Value = 'Hiển thị thông báo lỗi và số lượng không thay đổi'
$ws2.Cells.Item(13, 7).Value = 'vượt số lượng tồn kho'
$ws2.Cells.Item(13, 8).Value = 'Fail'
$ws2.Cells.Item(14, 1).Value = 'TC13'
$ws2.Cells.Item(14, 2).Value = 'Giỏ hàng'
$ws2.Cells.Item(14, 3).Value = 'Chức năng'
$ws2.Cells.Item(14, 4).Value = 'Kiểm tra chức năng button "Thanh toán"'
$ws2.Cells.Item(14, 6).Value = 'Chuyển hướng đến trang thanh toán'
$ws2.Cells.Item(14, 7).Value = 'Chuyển hướng đến trang thanh toán'
$ws2.Cells.Item(14, 8).Value = 'Pass'
$ws2.Cells.Item(15, 1).Value = 'TC14'
$ws2.Cells.Item(15, 2).Value = 'Giỏ hàng'
$ws2.Cells.Item(15, 3).Value = 'Chức năng'
$ws2.Cells.Item(15, 4).Value = 'Kiểm tra chức năng liên kết "Tiếp tục mua sắm"'
$ws2.Cells.Item(15, 6).Value = 'Chuyển hướng đến trang sản phẩm'
$ws2.Cells.Item(15, 7).Value = 'Chuyển hướng đến trang sản phẩm'
$ws2.Cells.Item(15, 8).Value = 'Pass'
$ws2.Cells.Item(16, 1).Value = 'TC15'
$ws2.Cells.Item(16, 2).Value = 'Giỏ hàng'
$ws2.Cells.Item(16, 3).Value = 'Chức năng'
$ws2.Cells.Item(16, 4).Value = 'Kiểm tra chức năng button "Xóa"'
$ws2.Cells.Item(16, 6).Value = 'Sản phẩm được xóa khỏi giỏ hàng và cập nhật tổng tiền'
$ws2.Cells.Item(16, 7).Value = 'Sản phẩm bị xóa khỏi giỏ hàng thành công'
$ws2.Cells.Item(16, 8).Value = 'Fail'
$ws2.Cells.Item(17, 1).Value = 'TC16'
$ws2.Cells.Item(17, 2).Value = 'Giỏ hàng'
$ws2.Cells.Item(17, 3).Value = 'Chức năng'
$ws2.Cells.Item(17, 4).Value = 'Kiểm tra xóa sản phẩm cuối cùng trong giỏ hàng'
$ws2.Cells.Item(17, 6).Value = 'Giỏ hàng trống, hiển thị thông báo (nếu có)'
$ws2.Cells.Item(17, 7).Value = 'Sản phẩm bị xóa khỏi giỏ hàng thành công'
$ws2.Cells.Item(17, 8).Value = 'Fail'
$ws2.Cells.Item(18, 1).Value = 'TC17'
$ws2.Cells.Item(18, 2).Value = 'Giỏ hàng'
$ws2.Cells.Item(18, 3).Value = 'Chức năng'
$ws2.Cells.Item(18, 4).Value = 'Kiểm tra tổng tiền khi giỏ hàng trống'
$ws2.Cells.Item(18, 6).Value = 'Tổng tiền hiển thị là 0'
$ws2.Cells.Item(18, 7).Value = 'Giỏ hàng không trống'
$ws2.Cells.Item(18, 8).Value = 'Fail'

# Make "Demo" the active sheet/tab (matches activeTab="5" in the target workbook).
$ws1.Activate()
